# Actualización automática de datos
# Adds a new transaction row (row 10) to the register sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 45707
$ws.Range("A10").NumberFormat = "yyyy-mm-dd"

$ws.Range("C10").Value = "Tinte"
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = "viriginia"
$ws.Range("G10").Value = "Efectivo"
